$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the already-styled column-A cell (bold, centered,
# bordered) down onto the two new rows (A4, A5) before writing their values.
$ws.Range("A2").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A (group id) values
$ws.Range("A2").Value = 3
$ws.Range("A3").Value = 0
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 2

# Column B (count) values
$ws.Range("B2").Value = 241
$ws.Range("B3").Value = 194
$ws.Range("B4").Value = 191
$ws.Range("B5").Value = 89
